$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet / "as-of" date references from 2021-11-17 to 2021-11-18
$ws.Name = "Through 2021-11-18"
$ws.Range("A12").Value = "November (through 11-18)"

# Update the November row (row 12) with the new daily totals
$ws.Range("B12").Value = 18
$ws.Range("D12").Value = 74
$ws.Range("E12").Value = 39
$ws.Range("F12").Value = 28
$ws.Range("G12").Value = 109
$ws.Range("H12").Value = 125

# Update the Total row (row 13) to reflect the new cumulative totals
$ws.Range("B13").Value = 276
$ws.Range("D13").Value = 784
$ws.Range("E13").Value = 654
$ws.Range("F13").Value = 510
$ws.Range("G13").Value = 1166
$ws.Range("H13").Value = 1567
